$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the twelve-month period headers: drop 1396/12, add 1401/12 (rolling one year forward)
$ws.Range("E8").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F8").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G8").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H8").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I8").Value = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E17").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F17").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G17").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H17").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I17").Value = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E27").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F27").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G27").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H27").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I27").Value = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E36").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F36").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G36").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H36").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I36").Value = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E44").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F44").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G44").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H44").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I44").Value = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E53").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F53").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G53").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H53").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I53").Value = "دوازده ماهه منتهی به 1401/12"

# Update the unit label for the sell-rate table ("kg/rial" -> "ton/rial") for both rows
$ws.Range("C39").Value = "تن / ریال"
$ws.Range("C40").Value = "تن / ریال"

# Shift each data row one period to the left and fill in the new 1401/12 figures
# Row 10
$ws.Range("E10").Value = "-"
$ws.Range("F10").Value = "-"
$ws.Range("G10").Value = "-"
$ws.Range("H10").Value = 49199
$ws.Range("I10").Value = 56894

# Row 11
$ws.Range("E11").Value = 38071299
$ws.Range("F11").Value = 36735491
$ws.Range("G11").Value = 41413
$ws.Range("H11").Value = "-"
$ws.Range("I11").Value = "-"

# Row 12
$ws.Range("E12").Value = "-"
$ws.Range("F12").Value = "-"
$ws.Range("G12").Value = "-"
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0

# Row 13
$ws.Range("E13").Value = 38071299
$ws.Range("F13").Value = 36735491
$ws.Range("G13").Value = 41413
$ws.Range("H13").Value = 49199
$ws.Range("I13").Value = 56894

# Row 19
$ws.Range("E19").Value = -17240
$ws.Range("F19").Value = -2750
$ws.Range("G19").Value = -5000
$ws.Range("H19").Value = "-"
$ws.Range("I19").Value = "-"

# Row 20
$ws.Range("E20").Value = "-"
$ws.Range("F20").Value = "-"
$ws.Range("G20").Value = "-"
$ws.Range("H20").Value = 43134
$ws.Range("I20").Value = 60405

# Row 21
$ws.Range("E21").Value = 37896395
$ws.Range("F21").Value = 35400251
$ws.Range("G21").Value = 43283119
$ws.Range("H21").Value = "-"
$ws.Range("I21").Value = "-"

# Row 22
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = -55

# Row 23
$ws.Range("E23").Value = 37879155
$ws.Range("F23").Value = 35397501
$ws.Range("G23").Value = 43278119
$ws.Range("H23").Value = 43134
$ws.Range("I23").Value = 60350

# Row 29
$ws.Range("E29").Value = -1223
$ws.Range("F29").Value = -3748
$ws.Range("G29").Value = -525
$ws.Range("H29").Value = "-"
$ws.Range("I29").Value = "-"

# Row 30
$ws.Range("E30").Value = 2584862
$ws.Range("F30").Value = 3302464
$ws.Range("G30").Value = 5322034
$ws.Range("H30").Value = 10593475
$ws.Range("I30").Value = 20425052

# Row 31
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = -19267

# Row 32
$ws.Range("E32").Value = 2583639
$ws.Range("F32").Value = 3298716
$ws.Range("G32").Value = 5321509
$ws.Range("H32").Value = 10593475
$ws.Range("I32").Value = 20405785

# Row 38
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = "-"
$ws.Range("I38").Value = "-"

# Row 39
$ws.Range("E39").Value = 68209
$ws.Range("F39").Value = 93289
$ws.Range("G39").Value = 122959
$ws.Range("H39").Value = 245594543
$ws.Range("I39").Value = 338135121

# Row 40
$ws.Range("E40").Value = 0
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0

# Row 46
$ws.Range("E46").Value = 843
$ws.Range("F46").Value = 2928
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = "-"
$ws.Range("I46").Value = "-"

# Row 47
$ws.Range("E47").Value = -1780804
$ws.Range("F47").Value = -2580156
$ws.Range("G47").Value = -3884461
$ws.Range("H47").Value = -8499119
$ws.Range("I47").Value = -16273879

# Row 48
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0

# Row 49
$ws.Range("E49").Value = -1779961
$ws.Range("F49").Value = -2577228
$ws.Range("G49").Value = -3884461
$ws.Range("H49").Value = -8499119
$ws.Range("I49").Value = -16273879

# Row 55
$ws.Range("E55").Value = -380
$ws.Range("F55").Value = -820
$ws.Range("G55").Value = -525
$ws.Range("H55").Value = "-"
$ws.Range("I55").Value = "-"

# Row 56
$ws.Range("E56").Value = 804058
$ws.Range("F56").Value = 722308
$ws.Range("G56").Value = 1437573
$ws.Range("H56").Value = 2094356
$ws.Range("I56").Value = 4151173

# Row 57
$ws.Range("E57").Value = "-"
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = -19267

# Row 58
$ws.Range("E58").Value = 803678
$ws.Range("F58").Value = 721488
$ws.Range("G58").Value = 1437048
$ws.Range("H58").Value = 2094356
$ws.Range("I58").Value = 4131906
